$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: extend sequence 0..9 (B1:K1) to 0..11 (B1:M1) ---
$ws.Range("K1").Copy()
$ws.Range("L1:M1").PasteSpecial(-4122)
$ws.Range("L1").Value = 10
$ws.Range("M1").Value = 11

# --- Row 2 (course 107) ---
$ws.Range("G2").Value = "Practicum"
$ws.Range("H2").ClearContents()
$ws.Range("I2").Value = "in"
$ws.Range("J2").Value = "Journalism"
$ws.Range("K2").Value = "TBA"
$ws.Range("L2").Value = "TBA"
$ws.Range("M2").Value = "Cooper Jan"

# --- Row 3 (course 308) ---
$ws.Range("G3").Value = "Writing"
$ws.Range("H3").Value = "Fellowship"
$ws.Range("I3").Value = "TR"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "0300"
$ws.Range("J3").Style = "Normal"
$ws.Range("K3").Value = "0350pm"
$ws.Range("L3").Value = "CARN"
$ws.Range("M3").Value = "Petzak Nicholas"
